# Fruta / hortaliza, semanal
# Insert a new weekly record as row 54, pushing existing rows 54:78 down to 55:79.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 54 (shifts rows 54-78 down to 55-79,
# and copies the date-format style from row 54 into the new row's D cell).
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new weekly record.
$ws.Range("A54").Value = 2
$ws.Range("B54").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C54").Value = "Coquimbo"
$ws.Range("D54").Value = 44839
$ws.Range("E54").Value = 4
$ws.Range("F54").Value = 100112022
$ws.Range("G54").Value = "Arveja Verde"
$ws.Range("H54").Value = "Perfection"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 700
$ws.Range("K54").Value = 22000
$ws.Range("L54").Value = 24000
$ws.Range("M54").Value = 23000
$ws.Range("N54").Value = "$/malla 25 kilos"
$ws.Range("O54").Value = "Provincia de Limarí"
$ws.Range("P54").Value = 920
$ws.Range("Q54").Value = 25
$ws.Range("R54").Value = "Hortaliza"
